# Temporary fix for non-display of tweets
# Populates the negative-tweets table (rows 4-16) and turns each "url"
# cell into a clickable hyperlink (underlined, blue) pointing at itself.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=4;  A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/905906972393963520"; E="#0daytoday #Apache Struts 2 REST Plugin XStream Remote Code Execution #Exploit #0day https://t.co/rwaxBrUvuC"},
    @{Row=5;  A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/905906732098088966"; E="#0daytoday #SourceTree Remote Code Execution #Exploit https://t.co/cGo5U8zYGm"},
    @{Row=6;  A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/905878605128884225"; E="#0daytoday #Gh0st Client - Buffer Overflow Exploit [remote #exploits #0day #Exploit] https://t.co/ce22CcFpwa"},
    @{Row=7;  A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/905878602645897217"; E="#0daytoday #PlugX Controller Stack Overflow Exploit [remote #exploits #0day #Exploit] https://t.co/jcciGj0kBX"},
    @{Row=8;  A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/905750181722947584"; E="#0daytoday #Apache #Struts 2.5 - Remote Code Execution #0day #Exploit https://t.co/NRIvcrWlEq"},
    @{Row=9;  A="Sep 09 2017"; B="negative"; C=-0.05; D="https://twitter.com/statuses/905728881197408257"; E="#0daytoday #WordPress Gym Management System Code Execution / Remote Cross Site Scripting Vulnerabil [#0day #Exploit] https://t.co/6BbjaqVOyR"},
    @{Row=10; A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/905513948224438272"; E="#0daytoday #Tor - Linux Sandbox Breakout via X11 Exploit [remote #exploits #0day #Exploit] https://t.co/8bCWbuT4Uv"},
    @{Row=11; A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/905505647285813252"; E="#0daytoday #Jungo DriverWizard WinDriver - Kernel Pool Overflow Exploit [remote #exploits #0day #Exploit] https://t.co/jc9QMlrawK"},
    @{Row=12; A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/905355158887907328"; E="#0daytoday #Samsung Internet Browser - SOP Bypass Exploit [remote #exploits #0day #Exploit] https://t.co/ka6tQuu1ZI"},
    @{Row=13; A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/905079577915088897"; E="#0daytoday #NEC EXPRESS CLUSTER clpwebmc Remote #Root #Exploit #0day https://t.co/4i3pB8Hwp8"},
    @{Row=14; A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/904803563603406850"; E="#0daytoday #RubyGems < 2.6.13 - Arbitrary File Overwrite #Exploit https://t.co/sM5XxS3Rxt"},
    @{Row=15; A="Sep 09 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/904743627418513408"; E="#0daytoday #Wireless Repeater BE126 - Remote Code Execution Exploit [webapps #exploits #0day #Exploit] https://t.co/MJBKKmDAxw"},
    @{Row=16; A="Aug 08 2017"; B="negative"; C=-0.1;  D="https://twitter.com/statuses/902990838938058752"; E="#0daytoday #Joomla Joomanager 2.0.0 Component - Arbitrary File Download Vulnerability [#0day #Exploit] https://t.co/OJCfpL4p9I"}
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E

    $urlCell = $ws.Range("D$row")
    $ws.Hyperlinks.Add($urlCell, $r.D)
    $urlCell.Font.Underline = 2
    $urlCell.Font.Color = 16711680
}
